# Added periodic & upfront related scenarios
# The "repaymentstrategy" sample value on the ProductLoanInput sheet is
# changed from "Mifos style" to "Penalties, Fees, Interest, Principal order",
# and the cell is given left/top alignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131   # xlLeft
$cell.VerticalAlignment = -4160     # xlTop

# Leave the selection on the edited cell, as Excel would after the edit.
$ws.Activate()
$cell.Select()
